$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Funcionário"
$ws.Range("D9").Value = "3. Mostra lista de carros em produção ou espera"

$ws.Range("D10").Select()
